$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the title paragraph currently reads
#   "Lenth's Critical Value" + "에 대한 이론, 수식, 사례"
# (the Korean tail is built from several separate runs). The edit drops
# everything after "Lenth's Critical Value", leaving only that run.
# ---------------------------------------------------------------------
$title = "Lenth's Critical Value"
$p1 = $d.Paragraphs(1)
$full1 = $p1.Range
if ($full1.Text.StartsWith($title)) {
    # Range covering everything after the title text up to (but not
    # including) the paragraph mark at the very end of the paragraph.
    $tail = $d.Range($full1.Start + $title.Length, $full1.End - 1)
    if ($tail.Start -lt $tail.End) {
        $tail.Delete()
    }
}

# ---------------------------------------------------------------------
# Change 2: the paragraph "수식:<line break>Lenth's Critical Value를 ..."
# currently holds "Lenth's Critical Value" as a single run (right after
# the manual line break). The edit splits it into two runs:
#   "Lenth's " (keeps the leading <w:br/>) and "Critical Value".
# Locate the paragraph unambiguously by matching the manual line break
# (vertical-tab, chr 11) immediately followed by the phrase.
# ---------------------------------------------------------------------
$lead = "Lenth's "
$rest = "Critical Value"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    $brIdx = $t.IndexOf([char]11)
    if ($brIdx -ge 0) {
        $after = $t.Substring($brIdx + 1)
        if ($after.StartsWith($lead + $rest)) {
            $full = $p.Range
            $start = $full.Start + $brIdx + 1 + $lead.Length
            $splitRange = $d.Range($start, $start + $rest.Length)
            # Force the run boundary to appear exactly at the
            # "Lenth's " / "Critical Value" split point: flipping a
            # character property on just this sub-range and then
            # flipping it straight back splits the run (as a side
            # effect of the assignment) without altering any visible
            # formatting, since the net property value is unchanged.
            $orig = $splitRange.Bold
            $splitRange.Bold = 1
            $splitRange.Bold = $orig
            break
        }
    }
}
